$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new text value, and whether the value must be
# forced to Text format first (values that Excel would otherwise parse as
# a number/date, e.g. "0.996" or "1.00") so it round-trips as a string,
# matching the original inlineStr cell - exactly like the source diff.
$updates = @(
    @{ Cell = "D2"; Value = "66.728.63"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -1.90%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.212.36"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -3.49%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "0.996"; ForceText = $true },
    @{ Cell = "E4"; Value = "  -0.31%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "579.56"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -3.94%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "140.13"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -13.98%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.997"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -0.33%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "3.205.91"; ForceText = $false },
    @{ Cell = "E8"; Value = "  -3.48%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.522"; ForceText = $true },
    @{ Cell = "E9"; Value = "  -9.04%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.161"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -11.87%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "6.42"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -2.59%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "0.476"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -9.58%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.0000231"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -9.47%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "35.53"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -14.41%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "3.724.03"; ForceText = $false },
    @{ Cell = "E15"; Value = "  -3.68%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "66.618.90"; ForceText = $false },
    @{ Cell = "E16"; Value = "  -2.19%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "3.200.09"; ForceText = $false },
    @{ Cell = "E17"; Value = "  -3.82%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "0.112"; ForceText = $true },
    @{ Cell = "E18"; Value = "  -4.41%  "; ForceText = $false },
    @{ Cell = "B19"; Value = "BitcoinCash"; ForceText = $false },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; ForceText = $false },
    @{ Cell = "D19"; Value = "501.56"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -8.85%  "; ForceText = $false },
    @{ Cell = "B20"; Value = "Polkadot"; ForceText = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; ForceText = $false },
    @{ Cell = "D20"; Value = "6.76"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -11.86%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "14.16"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -11.68%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "0.711"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -10.94%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "7.31"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -12.29%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "81.61"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -8.87%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "12.76"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -10.52%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E26"; Value = "  +0.23%  "; ForceText = $false },
    @{ Cell = "E27"; Value = "  -11.85%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "27.63"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -10.22%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "2.02"; ForceText = $true },
    @{ Cell = "E29"; Value = "  -10.77%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "7.47"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -7.97%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "1.17"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -2.17%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "2.48"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -6.64%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -0.47%  "; ForceText = $false },
    @{ Cell = "B34"; Value = "Filecoin"; ForceText = $false },
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; ForceText = $false },
    @{ Cell = "D34"; Value = "6.05"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -16.77%  "; ForceText = $false },
    @{ Cell = "B35"; Value = "Bittensor"; ForceText = $false },
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; ForceText = $false },
    @{ Cell = "D35"; Value = "492.26"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -14.28%  "; ForceText = $false },
    @{ Cell = "B36"; Value = "OKB"; ForceText = $false },
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; ForceText = $false },
    @{ Cell = "D36"; Value = "53.86"; ForceText = $true },
    @{ Cell = "E36"; Value = "  -2.44%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "5.26"; ForceText = $true },
    @{ Cell = "E37"; Value = "  -14.86%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.0413"; ForceText = $true },
    @{ Cell = "E38"; Value = "  -7.70%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "0.0808"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -10.94%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "8.44"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -14.31%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "0.118"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -13.44%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "2.818.99"; ForceText = $false },
    @{ Cell = "E42"; Value = "  -7.42%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "2.53"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -11.42%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "0.249"; ForceText = $true },
    @{ Cell = "E45"; Value = "  -10.05%  "; ForceText = $false },
    @{ Cell = "B46"; Value = "InjectiveProtocol"; ForceText = $false },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; ForceText = $false },
    @{ Cell = "D46"; Value = "24.82"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -13.78%  "; ForceText = $false },
    @{ Cell = "B47"; Value = "Fetch.AI"; ForceText = $false },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; ForceText = $false },
    @{ Cell = "D47"; Value = "2.02"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -10.05%  "; ForceText = $false },
    @{ Cell = "B48"; Value = "Monero"; ForceText = $false },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; ForceText = $false },
    @{ Cell = "D48"; Value = "120.63"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -5.54%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "0.0₃0526"; ForceText = $false },
    @{ Cell = "E49"; Value = "  -15.19%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "0.108"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -9.88%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "2.12"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -20.05%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Temporarily set Text number format so the numeric-looking string is
        # not auto-converted into a number, then restore the default style so
        # no stray format/style is left behind on the cell.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
